# Insert a new observation row at row 2 (pushing existing rows 2-5 down to 3-6)
# and populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

function Set-TextCell {
    # Writes a literal text value into a cell without letting Excel's
    # "smart" typing heuristics reinterpret numeric- or date-looking
    # strings as numbers/dates (which would also pull in a new number
    # format style). We build a text formula ="..." and flatten it to a
    # plain value via copy/paste-special, which keeps the cell's style
    # untouched (stays on the shared default "Normal" style).
    param(
        [string]$addr,
        [string]$text
    )
    $escaped = $text -replace '"', '""'
    $r = $ws.Range($addr)
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)  # xlPasteValues
}

# Numeric cells
$ws.Range("A2").Value = 55476489
$ws.Range("B2").Value = 96309
$ws.Range("E2").Value = 219797
$ws.Range("Q2").Value = 610865.2120036986
$ws.Range("R2").Value = 6566445.901323199
$ws.Range("S2").Value = 10

# Text cells
Set-TextCell "C2" "Ovaliderad"
Set-TextCell "D2" "LC"
Set-TextCell "F2" "Purpurknipprot"
Set-TextCell "G2" "Epipactis atrorubens"
Set-TextCell "H2" "(Hoffm.) Besser"
Set-TextCell "I2" "1"
Set-TextCell "J2" "stjälkar/strån/skott"
Set-TextCell "K2" "blomknopp"
Set-TextCell "P2" "Magsjötorp, N om, Srm"
Set-TextCell "T2" "Södermanland"
Set-TextCell "U2" "Strängnäs"
Set-TextCell "V2" "Södermanland"
Set-TextCell "W2" "Länna"
Set-TextCell "Y2" "2013-06-29"
Set-TextCell "Z2" "00:00"
Set-TextCell "AA2" "2013-06-29"
Set-TextCell "AB2" "00:00"
Set-TextCell "AI2" "Tallskog, gles"
Set-TextCell "AW2" "Håkan Gustafson"
Set-TextCell "AX2" "Håkan Gustafson"

# Boolean cells
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
